{"js": "const body = context.document.body;\n\n// The lab-report number run currently reads \"2\" \u2014 bump it to \"3\".\n// matchWholeWord avoids also matching the \"2\" inside \"2018\" elsewhere\n// in the document.\nconst numberMatches = body.search(\"2\", { matchWholeWord: true });\nnumberMatches.load(\"items\");\nawait context.sync();\n\nconst numberRun = numberMatches.items[0];\nnumberRun.insertText(\"3\", \"Replace\");\nawait context.sync();\n\n// Word's \"_GoBack\" bookmark always marks the location of the most\n// recent edit, so after changing the number it moves from the end of\n// the document (after \"\u041e\u0440\u0451\u043b, 2018\") to right after the new \"3\".\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst newNumberMatches = body.search(\"3\", { matchWholeWord: true });\nnewNumberMatches.load(\"items\");\nawait context.sync();\n\nconst newNumberRun = newNumberMatches.items[0];\nconst collapsedEnd = newNumberRun.getRange(\"End\");\ncollapsedEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Step 1: bump the lab-report number from \"2\" to \"3\" -------------------\n# MatchWholeWord keeps this from also hitting the \"2\" inside \"2018\" later\n# in the document.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"2\"\n$find.MatchWholeWord = $true\n$find.Execute() | Out-Null\n$rng.Text = \"3\"\n\n# --- Step 2: move the \"_GoBack\" bookmark to sit right after the new \"3\" ---\n# Word always keeps \"_GoBack\" collapsed at the location of the most recent\n# edit, so after changing the number it has to move from the end of the\n# document (after \"\u041e\u0440\u0451\u043b, 2018\") to immediately after the new \"3\".\n#\n# $rng is now collapsed right after \"3\". We temporarily insert a sentinel\n# character after it so the target bookmark position is not the very last\n# character of the paragraph (i.e. not immediately before the paragraph\n# mark) while we create the bookmark, then remove the sentinel afterwards.\n# The bookmark stays anchored to the correct spot once the sentinel goes\n# away.\n$rng.InsertAfter(\"X\")\n$bookmarkPos = $d.Range($rng.End - 1, $rng.End - 1)\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkPos)\n\n$bk = $d.Bookmarks(\"_GoBack\")\n$sentinel = $d.Range($bk.End, $bk.End + 1)\n$sentinel.Delete()\n"}
